{"js": "// Fix: Modified documentation\n// Replace \"allowing for\" with \"resulting in\" in the WishlistOperations\n// bullet describing the Dependency Inversion Principle.\nconst body = context.document.body;\nconst results = body.search(\"allowing for\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Text \"allowing for\" not found in document body.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"resulting in\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix : Modified documentation\n#\n# WishlistOperations bullet (Dependency Inversion Principle section):\n#   \"...WishlistServices class, allowing for switching of wishlist\n#   implementations.\"\n# becomes\n#   \"...WishlistServices class, resulting in switching of wishlist\n#   implementations.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"allowing for\"\n$find.Replacement.Text = \"resulting in\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\n# wdFindContinue = 0, wdReplaceAll = 2\n$found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, `\n  $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, `\n  $find.Replacement.Text, 2)\n\nif (-not $found) {\n  throw 'Text \"allowing for\" was not found in the document.'\n}\n"}
